$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.066.49'
$ws.Range("E2").Value = '  +4.23%  '

$ws.Range("D3").Value = '3.253.42'
$ws.Range("E3").Value = '  +2.49%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '''577.90'
$ws.Range("E5").Value = '  +2.33%  '

$ws.Range("D6").Value = '''178.16'
$ws.Range("E6").Value = '  +4.52%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("E8").Value = '  -1.05%  '

$ws.Range("D9").Value = '3.250.16'
$ws.Range("E9").Value = '  +2.52%  '

$ws.Range("E10").Value = '  +4.28%  '

$ws.Range("E11").Value = '  +1.75%  '

$ws.Range("D12").Value = '''0.413'
$ws.Range("E12").Value = '  +4.19%  '

$ws.Range("D13").Value = '3.816.39'
$ws.Range("E13").Value = '  +2.52%  '

$ws.Range("D14").Value = '''0.137'
$ws.Range("E14").Value = '  +0.65%  '

$ws.Range("D15").Value = '''28.16'
$ws.Range("E15").Value = '  +2.68%  '

$ws.Range("D16").Value = '67.035.11'
$ws.Range("E16").Value = '  +4.22%  '

$ws.Range("E17").Value = '  +2.84%  '

$ws.Range("D18").Value = '3.255.60'
$ws.Range("E18").Value = '  +2.68%  '

$ws.Range("E19").Value = '  +2.09%  '

$ws.Range("E20").Value = '  +3.14%  '

$ws.Range("D21").Value = '''373.40'
$ws.Range("E21").Value = '  +5.43%  '

$ws.Range("D22").Value = '''7.65'
$ws.Range("E22").Value = '  +5.97%  '

$ws.Range("E23").Value = '  -0.13%  '

$ws.Range("D24").Value = '''70.80'
$ws.Range("E24").Value = '  +2.37%  '

$ws.Range("E25").Value = '  +1.52%  '

$ws.Range("D26").Value = '3.396.81'
$ws.Range("E26").Value = '  +2.81%  '

$ws.Range("E27").Value = '  -0.15%  '

$ws.Range("D28").Value = '''9.92'
$ws.Range("E28").Value = '  +3.63%  '

$ws.Range("E29").Value = '  +2.00%  '

$ws.Range("E30").Value = '  +0.15%  '

$ws.Range("E31").Value = '  +4.24%  '

$ws.Range("D32").Value = '''5.64'
$ws.Range("E32").Value = '  +0.29%  '

$ws.Range("E33").Value = '  +2.31%  '

$ws.Range("E34").Value = '  -0.05%  '

$ws.Range("D35").Value = '''1.26'
$ws.Range("E35").Value = '  +4.63%  '

$ws.Range("E36").Value = '  +2.73%  '

$ws.Range("D37").Value = '''166.58'
$ws.Range("E37").Value = '  +7.18%  '

$ws.Range("E38").Value = '  +4.44%  '

$ws.Range("E39").Value = '  +5.24%  '

$ws.Range("E40").Value = '  +10.12%  '

$ws.Range("D41").Value = '''27.08'
$ws.Range("E41").Value = '  +4.65%  '

$ws.Range("D42").Value = '''2.60'
$ws.Range("E42").Value = '  +1.96%  '

$ws.Range("D43").Value = '2.760.95'
$ws.Range("E43").Value = '  +6.03%  '

$ws.Range("D44").Value = '''6.49'
$ws.Range("E44").Value = '  +8.17%  '

$ws.Range("B45").Value = 'Filecoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D45").Value = '''4.39'
$ws.Range("E45").Value = '  +4.93%  '

$ws.Range("B46").Value = 'Bittensor'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D46").Value = '''353.75'
$ws.Range("E46").Value = '  +9.33%  '

$ws.Range("D47").Value = '''25.46'
$ws.Range("E47").Value = '  +6.04%  '

$ws.Range("D48").Value = '''40.48'
$ws.Range("E48").Value = '  +1.94%  '

$ws.Range("D49").Value = '''0.0675'
$ws.Range("E49").Value = '  +2.35%  '

$ws.Range("D50").Value = '''0.0280'
$ws.Range("E50").Value = '  +3.26%  '

$ws.Range("E51").Value = '  +0.95%  '
